# Refresh the cryptos list (Coin / Link / Price / Volume(1h)) with the
# latest scrape, mirroring the "Updated cryptos list ... with GitHub
# Actions" automation commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving it as TEXT (many
# "price" strings such as "1.006" or "0.4844" look like numbers and Excel
# would otherwise silently coerce them). We briefly force a text number
# format for the assignment, then clear formats again so we don't leave a
# stray style index behind on a cell that previously had none.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row -> updated Coin, Link, Price, Volume(1h). $null means "leave as is".
$updates = @(
    @{ Row = 2;  B = $null;          C = $null;                                                            D = '29.493.37';   E = '  +0.84%  ' }
    @{ Row = 3;  B = $null;          C = $null;                                                            D = '1.922.52';    E = '  +1.37%  ' }
    @{ Row = 4;  B = $null;          C = $null;                                                            D = '1.006';       E = '  +0.53%  ' }
    @{ Row = 5;  B = $null;          C = $null;                                                            D = $null;         E = '  +0.80%  ' }
    @{ Row = 6;  B = $null;          C = $null;                                                            D = $null;         E = '  +0.47%  ' }
    @{ Row = 7;  B = $null;          C = $null;                                                            D = '0.4844';      E = '  +3.09%  ' }
    @{ Row = 8;  B = $null;          C = $null;                                                            D = '0.4090';      E = '  +1.71%  ' }
    @{ Row = 9;  B = $null;          C = $null;                                                            D = '0.08180';     E = '  +2.25%  ' }
    @{ Row = 10; B = $null;          C = $null;                                                            D = '1.027';       E = '  +3.40%  ' }
    @{ Row = 11; B = $null;          C = $null;                                                            D = $null;         E = '  +5.80%  ' }
    @{ Row = 12; B = $null;          C = $null;                                                            D = '1.911.59';    E = '  +3.91%  ' }
    @{ Row = 13; B = $null;          C = $null;                                                            D = '6.053';       E = '  +3.41%  ' }
    @{ Row = 14; B = $null;          C = $null;                                                            D = '7.234';       E = '  +2.78%  ' }
    @{ Row = 15; B = $null;          C = $null;                                                            D = '91.59';       E = '  +2.92%  ' }
    @{ Row = 16; B = 'TRON';         C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx';             D = '0.06765';     E = '  +2.15%  ' }
    @{ Row = 17; B = 'BinanceUSD';   C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd';      D = '1.006';       E = '  +0.45%  ' }
    @{ Row = 18; B = $null;          C = $null;                                                            D = '0.00001040';  E = '  +1.37%  ' }
    @{ Row = 19; B = $null;          C = $null;                                                            D = '17.81';       E = '  +1.87%  ' }
    @{ Row = 20; B = $null;          C = $null;                                                            D = $null;         E = '  +0.49%  ' }
    @{ Row = 21; B = $null;          C = $null;                                                            D = '29.520.44';   E = '  +0.88%  ' }
    @{ Row = 22; B = $null;          C = $null;                                                            D = $null;         E = '  +2.49%  ' }
    @{ Row = 23; B = $null;          C = $null;                                                            D = '11.78';       E = '  +2.08%  ' }
    @{ Row = 24; B = $null;          C = $null;                                                            D = $null;         E = '  -0.73%  ' }
    @{ Row = 25; B = $null;          C = $null;                                                            D = '2.134.97';    E = '  +0.23%  ' }
    @{ Row = 26; B = $null;          C = $null;                                                            D = '6.712';       E = '  +10.41%  ' }
    @{ Row = 27; B = $null;          C = $null;                                                            D = '156.71';      E = '  +1.37%  ' }
    @{ Row = 28; B = $null;          C = $null;                                                            D = '20.09';       E = '  +2.17%  ' }
    @{ Row = 29; B = $null;          C = $null;                                                            D = '2.128';       E = '  +2.11%  ' }
    @{ Row = 30; B = $null;          C = $null;                                                            D = '120.74';      E = '  +3.01%  ' }
    @{ Row = 31; B = $null;          C = $null;                                                            D = '1.027';       E = '  -2.74%  ' }
    @{ Row = 32; B = $null;          C = $null;                                                            D = '0.09567';     E = '  +1.38%  ' }
    @{ Row = 33; B = $null;          C = $null;                                                            D = $null;         E = '  +3.69%  ' }
    @{ Row = 34; B = $null;          C = $null;                                                            D = '3.573';       E = '  +0.98%  ' }
    @{ Row = 35; B = $null;          C = $null;                                                            D = $null;         E = '  -0.01%  ' }
    @{ Row = 36; B = $null;          C = $null;                                                            D = $null;         E = '  +2.12%  ' }
    @{ Row = 37; B = $null;          C = $null;                                                            D = '0.06151';     E = '  +1.28%  ' }
    @{ Row = 38; B = $null;          C = $null;                                                            D = '1.183';       E = '  +0.64%  ' }
    @{ Row = 39; B = $null;          C = $null;                                                            D = '0.5994';      E = '  +3.17%  ' }
    @{ Row = 40; B = $null;          C = $null;                                                            D = '10.85';       E = '  +8.16%  ' }
    @{ Row = 41; B = $null;          C = $null;                                                            D = '8.016';       E = '  -0.46%  ' }
    @{ Row = 42; B = $null;          C = $null;                                                            D = '0.1866';      E = '  +2.23%  ' }
    @{ Row = 43; B = $null;          C = $null;                                                            D = '2.444';       E = '  -1.44%  ' }
    @{ Row = 44; B = $null;          C = $null;                                                            D = $null;         E = '  +2.61%  ' }
    @{ Row = 45; B = $null;          C = $null;                                                            D = '0.07625';     E = '  -0.93%  ' }
    @{ Row = 46; B = $null;          C = $null;                                                            D = '12.45';       E = '  +2.17%  ' }
    @{ Row = 47; B = $null;          C = $null;                                                            D = '0.5597';      E = '  +2.29%  ' }
    @{ Row = 48; B = $null;          C = $null;                                                            D = '1.965';       E = '  +3.44%  ' }
    @{ Row = 49; B = $null;          C = $null;                                                            D = '116.88';      E = '  +3.14%  ' }
    @{ Row = 50; B = $null;          C = $null;                                                            D = '2.435';       E = '  +4.34%  ' }
    @{ Row = 51; B = $null;          C = $null;                                                            D = '72.88';       E = '  +2.88%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.B) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($null -ne $u.C) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($null -ne $u.D) { Set-TextValue $ws.Cells.Item($r, 4) $u.D }
    if ($null -ne $u.E) { Set-TextValue $ws.Cells.Item($r, 5) $u.E }
}
